$d = $word.ActiveDocument

# Paragraph 1: "Haematococcus" / "Pluvialis" heading - add spell-check proofErr marks
# and split the <w:br/> off into its own run.
$xmlP1 = @'
<w:p w14:paraId="6CEF61E9" w14:textId="0E6B213E" w:rsidR="00384785" w:rsidRDefault="00384785" w:rsidP="00384785"><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:t>Haematococcus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:t>Pluvialis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:br/><w:t xml:space="preserve">Organic material from </w:t></w:r><w:r w:rsidR="001B7B6F"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:t xml:space="preserve">a </w:t></w:r><w:r w:rsidR="001B7B6F" w:rsidRPr="001B7B6F"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:noProof/><w:lang w:val="en"/></w:rPr><w:t>pr</w:t></w:r><w:r w:rsidR="001B7B6F"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:noProof/><w:lang w:val="en"/></w:rPr><w:t>i</w:t></w:r><w:r w:rsidR="001B7B6F" w:rsidRPr="001B7B6F"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:noProof/><w:lang w:val="en"/></w:rPr><w:t>stine</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00D80E8E"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:noProof/><w:lang w:val="en"/></w:rPr><w:t>highland</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$d.Paragraphs.Item(1).Range.InsertXML($xmlP1)

# Paragraph 4: "We produce 100% natural astaxanthin ..." - add spell-check proofErr
# marks around astaxanthin/Heamatococcus/pluvialis/photobioreactors/technologies,
# split "CO2"/"technologies" runs, and drop the old _GoBack bookmark (it is
# relocated to the "1900m altitude" paragraph below).
$xmlP4 = @'
<w:p w14:paraId="617F0F0A" w14:textId="29A19BA8" w:rsidR="00523FDE" w:rsidRPr="00523FDE" w:rsidRDefault="00384785" w:rsidP="00523FDE"><w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en"/></w:rPr><w:br/></w:r><w:r><w:t xml:space="preserve">We produce 100% natural </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>astaxanthin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Heamatococcus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pluvialis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> with the highest quality standards throughout a holistic production line. Microalgae biomass is cultivated in closed tubular </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>photobioreactors</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> with highly purified water</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>and C</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>O2</w:t></w:r><w:r><w:t xml:space="preserve"> to guarantee a contamination-free environment. We are dedicated to </w:t></w:r><w:r w:rsidR="00523FDE" w:rsidRPr="001B7B6F"><w:rPr><w:noProof/></w:rPr><w:t>produc</w:t></w:r><w:r w:rsidR="001B7B6F"><w:rPr><w:noProof/></w:rPr><w:t>ing</w:t></w:r><w:r><w:t xml:space="preserve"> the world’s highest quality </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>astaxanthin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for human consumption with our cutting-edge tec</w:t></w:r><w:r><w:t>hnolog</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ies</w:t></w:r><w:r><w:t xml:space="preserve"> and rich experience in microalgae products.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(4).Range.InsertXML($xmlP4)

# Paragraph 6: "World's Highest Quality Astaxanthin Producer" - proofErr around Astaxanthin.
$xmlP6 = @'
<w:p w14:paraId="3957D24B" w14:textId="77777777" w:rsidR="00523FDE" w:rsidRPr="00CA19C1" w:rsidRDefault="00384785" w:rsidP="00523FDE"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:lang w:val="en"/></w:rPr><w:br/></w:r><w:r w:rsidR="00523FDE" w:rsidRPr="00CA19C1"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>World</w:t></w:r><w:r w:rsidR="00523FDE" w:rsidRPr="00CA19C1"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">’s Highest Quality </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Astaxanthin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Producer </w:t></w:r></w:p>
'@
$d.Paragraphs.Item(6).Range.InsertXML($xmlP6)

# Paragraph 7: "Our world class 88,000 m2 astaxanthin facility ..." - proofErr around astaxanthin.
$xmlP7 = @'
<w:p w14:paraId="792B6760" w14:textId="1FDCC786" w:rsidR="00523FDE" w:rsidRDefault="00523FDE" w:rsidP="00523FDE"><w:r><w:t>Our world class 88,000 m</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>astaxanthin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> facility is located in a pristine mountainous area with natural fresh air, </w:t></w:r><w:r w:rsidR="001B7B6F" w:rsidRPr="001B7B6F"><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve">clean </w:t></w:r><w:r w:rsidRPr="001B7B6F"><w:rPr><w:noProof/></w:rPr><w:t>water</w:t></w:r><w:r w:rsidR="001B7B6F"><w:rPr><w:noProof/></w:rPr><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> and sufficient sunshine. </w:t></w:r></w:p>
'@
$d.Paragraphs.Item(7).Range.InsertXML($xmlP7)

# Paragraph 9: "Over 10tons Heamatococcus pluvialis annual production" - proofErr around
# Heamatococcus/pluvialis.
$xmlP9 = @'
<w:p w14:paraId="0C15C41E" w14:textId="77777777" w:rsidR="00523FDE" w:rsidRDefault="00523FDE" w:rsidP="00523FDE"><w:r><w:t xml:space="preserve">Over 10tons </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00BE1E58"><w:rPr><w:i/></w:rPr><w:t>Heamatococcus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00BE1E58"><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00BE1E58"><w:rPr><w:i/></w:rPr><w:t>pluvialis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> annual production </w:t></w:r></w:p>
'@
$d.Paragraphs.Item(9).Range.InsertXML($xmlP9)

# Paragraph 10: "Up to 7% astaxanthin concentration in microalgae powder" - proofErr around astaxanthin.
$xmlP10 = @'
<w:p w14:paraId="3A1161CE" w14:textId="77777777" w:rsidR="00523FDE" w:rsidRDefault="00523FDE" w:rsidP="00523FDE"><w:r><w:t xml:space="preserve">Up to 7% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>astaxanthin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> concentration in microalgae powder</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(10).Range.InsertXML($xmlP10)

# Paragraph 16: "1900m altitude with fresh air and water" - relocate the _GoBack bookmark here.
$xmlP16 = @'
<w:p w14:paraId="1B884A97" w14:textId="77777777" w:rsidR="00523FDE" w:rsidRDefault="00523FDE" w:rsidP="00523FDE"><w:r><w:t xml:space="preserve">1900m </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t>altitude with fresh air and water</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>
'@
$d.Paragraphs.Item(16).Range.InsertXML($xmlP16)
